$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lead")

# --- Copy row 2's cell formatting onto the new row 3 first (styles + row
# height), then fill values - so the quote-prefixed text cells land on an
# already-matching style instead of Excel minting a throwaway blended one.
$ws.Range("A2:U2").Copy()
$ws.Range("A3:U3").PasteSpecial(-4122)
$ws.Rows.Item(3).RowHeight = 52.8

# --- Fill the new row 3 values (identical to row 2 except the lead name) ---
# Use .Formula with a leading apostrophe for the numeric-looking text values
# ("001", "1245", "1", "0") so they stay text (shared strings) instead of
# collapsing into plain numbers, matching row 2's stored types.
$ws.Range("A3").Formula = "Active"
$ws.Range("B3").Formula = "Google"
$ws.Range("C3").Formula = "Admin Anh Tester"
$ws.Range("D3").Formula = "JSC_NEW"
$ws.Range("E3").Formula = "Yến Nhi 2"
$ws.Range("F3").Formula = "Đại Linh"
$ws.Range("G3").Formula = "Tester"
$ws.Range("H3").Formula = "Việt Nam"
$ws.Range("I3").Formula = "ngocnhi"
$ws.Range("J3").Formula = "Hà Nội"
$ws.Range("K3").Formula = "htester.com.vn"
$ws.Range("L3").Formula = "Vietnam"
$ws.Range("M3").Value = 982198605
$ws.Range("N3").Formula = "'001"
$ws.Range("O3").Formula = "'1245"
$ws.Range("P3").Formula = "Vietnamese"
$ws.Range("Q3").Formula = "NODO JSC"
$ws.Range("R3").Formula = "htest add new lead"
$ws.Range("S3").Value = 45971
$ws.Range("T3").Formula = "'1"
$ws.Range("U3").Formula = "'0"

# --- Make the Lead sheet the active tab/sheet and set its selection ---
$ws.Activate()
$ws.Range("G12").Select() | Out-Null
